$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold/border/centered, same as the other Row-1 headers)
# onto the three new header cells, then stamp in the header text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row (2-47).
# Every row on this roster shares the same team record.
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 78
    $ws.Cells.Item($row, 31).Value = 84
    $ws.Cells.Item($row, 32).Value = 0
}
